$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 14
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 12
$ws.Range("H2").Value = 5
